$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Supported Commands")
$ws.Columns.Item(4).ColumnWidth = 101.17
